# Backup QR Scanner data - 4/5/2025, 9:50:54 PM
#
# The QR-scanner app snapshots its in-memory scan log to a new worksheet
# each time a "backup" is triggered. This run adds one more sheet,
# "Dihdhdh", at the end of the workbook with the same 5-column layout
# (Number, Student ID, Location, Log Date, Log Time) used by every other
# backup sheet, holding the single scan that was just logged.

$wb = $excel.ActiveWorkbook

# New sheet goes after the current last sheet, so it ends up last/active.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Dihdhdh"

# Header row.
$ws.Range("A1").Value = "Number"
$ws.Range("B1").Value = "Student ID"
$ws.Range("C1").Value = "Location"
$ws.Range("D1").Value = "Log Date"
$ws.Range("E1").Value = "Log Time"

# Data row. "Number" is a real numeric sequence counter; the Student ID,
# Log Date and Log Time columns are stored as plain text in the source
# data (they can contain leading zeros / non-numeric IDs), so mark those
# cells as Text before writing them to stop Excel from reinterpreting a
# numeric-looking ID or a yyyy-mm-dd string as a number/date.
$ws.Range("A2").Value = 1

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "231249"

$ws.Range("C2").Value = "Dihdhdh"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2025-04-05"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "21:50:52"

$ws.Range("A1").Select()
